$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.963.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +10.57%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.529.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +14.01%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "490.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +17.97%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +23.58%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +14.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.520.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +14.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.37%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +16.64%  "

# Row 13
$ws.Range("E13").Value = "  +3.21%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.964.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +14.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.951.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +15.98%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +24.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.522.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +18.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +14.26%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +20.24%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.35%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +25.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.414"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +19.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.624.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0814"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +27.38%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +20.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +17.91%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.899"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.80%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.83%  "

# Row 38
$ws.Range("E38").Value = "  +20.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.619"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.61%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0561"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.45%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.75%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.020.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.84%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +30.41%  "

# Row 47
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "259.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +50.49%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0916"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +15.91%  "

# Row 49
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +21.70%  "
